$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: 8 rows (ECs/FAPs/MuSCs/Resolving-Mac) x (FAPs/Resolving-Mac target)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average expression value,
# H Ligand total expression value, I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value, K Receptor-expressing cells,
# L Receptor detection rate, M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value, P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity

$data = @(
    @("ECs","Il23a","Il23r","FAPs",3,1,3.109891,9.329673,0.4376734453228661,0.4376734453228661,1,0.3333333333333333,0.014443,0.043329,0.7689945869198687,0.7689945869198687,0.04491615571299999,0.404245401417,0.3365685102918531,0.3365685102918531),
    @("ECs","Il23a","Il23r","Resolving-Mac",3,1,3.109891,9.329673,0.4376734453228661,0.4376734453228661,1,0.3333333333333333,0.004338666666666667,0.013016,0.2310054130801314,0.2310054130801313,0.01349278041866667,0.121435023768,0.101104935031013,0.1011049350310129),
    @("FAPs","Il23a","Il23r","FAPs",3,1,2.963940333333333,8.891821,0.417132940486147,0.417132940486147,1,0.3333333333333333,0.014443,0.043329,0.7689945869198687,0.7689945869198687,0.04280819023433333,0.385273712109,0.3207729732598148,0.3207729732598148),
    @("FAPs","Il23a","Il23r","Resolving-Mac",3,1,2.963940333333333,8.891821,0.417132940486147,0.417132940486147,1,0.3333333333333333,0.004338666666666667,0.013016,0.2310054130801314,0.2310054130801313,0.01285954912622222,0.115735942136,0.09635996722633225,0.09635996722633224),
    @("MuSCs","Il23a","Il23r","FAPs",3,1,0.8168863333333333,2.450659,0.1149652691837634,0.1149652691837634,1,0.3333333333333333,0.014443,0.043329,0.7689945869198687,0.7689945869198687,0.01179828931233333,0.106184603811,0.08840766968609966,0.08840766968609966),
    @("MuSCs","Il23a","Il23r","Resolving-Mac",3,1,0.8168863333333333,2.450659,0.1149652691837634,0.1149652691837634,1,0.3333333333333333,0.004338666666666667,0.013016,0.2310054130801314,0.2310054130801313,0.003544197504888889,0.031897777544,0.02655759949766377,0.02655759949766376),
    @("Resolving-Mac","Il23a","Il23r","FAPs",2,0.6666666666666666,0.2147876666666667,0.644363,0.03022834500722351,0.03022834500722351,1,0.3333333333333333,0.014443,0.043329,0.7689945869198687,0.7689945869198687,0.003102178269666667,0.027919604427,0.02324543368210112,0.02324543368210111),
    @("Resolving-Mac","Il23a","Il23r","Resolving-Mac",2,0.6666666666666666,0.2147876666666667,0.644363,0.03022834500722351,0.03022834500722351,1,0.3333333333333333,0.004338666666666667,0.013016,0.2310054130801314,0.2310054130801313,0.0009318920897777779,0.008387028808,0.006982911325122394,0.006982911325122392)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}
